# Update the "Förändrad" (Changed) date column (C) from 2023-09-17 (45186)
# to 2023-09-19 (45188) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$oldSerial = 45186
$newSerial = 45188

# Data starts at row 2 (row 1 is the header) and the sheet's used range
# extends to row 218, matching the workbook's dimension.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
